$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the bug status for row 8 from "Unresolved" to "Resolved"
$ws.Range("E8").Value = "Resolved"

# Update the bug solution text for row 8
$ws.Range("F8").Value = "Layer masking wasn't set to ignore the IgnoreRaycast layer fixed"

# Update the active selection to F8
$ws.Range("F8").Select()
